$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8 is Oxalic acid. Add price data (USD/MT) and compute price/kg via formula,
# plus Location, updated Date, and a new Source link.
$ws.Range("M8").Value = 720
$ws.Range("N8").Formula = "=M8*0.88*0.001"
$ws.Range("P8").Value = "Germany"
$ws.Range("P8").NumberFormat = "0.0000"
$ws.Range("Q8").Value = "Q2 2024"
$ws.Range("R8").Value = "https://www.chemanalyst.com/Pricing-data/oxalic-acid-1556"

# Update the selected/active cell in the sheet view
$ws.Range("K16").Select()
